$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set priority (column C) values for the newly prioritized use cases (rows 1-3 per commit message)
$ws.Range("C6").Value = 4
$ws.Range("C9").Value = 4
$ws.Range("C10").Value = 4
$ws.Range("C11").Value = 4
$ws.Range("C12").Value = 4
$ws.Range("C13").Value = 4

# Move the view / selection to match where the author left off working
$ws.Application.ActiveWindow.ScrollRow = 37
$ws.Range("B49").Select()
